$wb = $excel.ActiveWorkbook

$data = @(
    @{ Sheet = "Neodymium";   C1 = 2030; C2 = 0.0002195405251500087; C3 = 0.01062411525673284;  C4 = 0.009608716352691784; C5 = 0.0000002138791829054013 },
    @{ Sheet = "Dysprosium";  C1 = 2030; C2 = 0.0002195405251500235; C3 = 0.01062411525673355;  C4 = 0.009608716352692431; C5 = 0.0000002138791829054185 },
    @{ Sheet = "Copper";      C1 = 2030; C2 = 0.007632681444695514;  C3 = 0.027535891297259;    C4 = 0.007370778697872926; C5 = 0.016160587324431 },
    @{ Sheet = "Raw silicon"; C1 = 2030; C2 = 0.005750015024097243;  C3 = 0.01921210602835477;  C4 = 0.0053924808017845;   C5 = 0.006847896595910315 }
)

foreach ($entry in $data) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $ws.Range("C1").Value = $entry.C1
    $ws.Range("C2").Value = $entry.C2
    $ws.Range("C3").Value = $entry.C3
    $ws.Range("C4").Value = $entry.C4
    $ws.Range("C5").Value = $entry.C5
}
